$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the two new journal entries (rows 29 & 30) -----------------------
# Set the values/formulas first (before touching formatting) so the I2
# SUM formula recalculates against the final numbers instead of a stale
# dependency snapshot taken mid-paste.

# Row 29: "Documentation" entry
$ws.Cells.Item(29, 1).Value = 45076
$ws.Cells.Item(29, 2).Formula = "=WEEKNUM(A29)"
$ws.Cells.Item(29, 3).Value = 2
$ws.Cells.Item(29, 4).Value = "Documentation"
$ws.Cells.Item(29, 5).Value = "documentation de toutes les fonctions"
$ws.Cells.Item(29, 6).Value = "aucun problème"

# Row 30: "Coding/implementation" entry
$ws.Cells.Item(30, 1).Value = 45076
$ws.Cells.Item(30, 2).Formula = "=WEEKNUM(A30)"
$ws.Cells.Item(30, 3).Value = 3.25
$ws.Cells.Item(30, 4).Value = "Coding/implementation"
$ws.Cells.Item(30, 5).Value = "tagging the code for the script"
$ws.Cells.Item(30, 6).Value = "aucun problème"

# Copy the formatting of the last existing data row (27) down onto the two
# new rows so they pick up the same styles (borders, alignment, number
# formats, etc.) used throughout the table.
$srcFormat = $ws.Range("A27:F27")
$dstFormat = $ws.Range("A29:F30")
$srcFormat.Copy()
$dstFormat.PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# --- Update the "Temps Total" sum to cover the new rows --------------------
$ws.Cells.Item(2, 9).Formula = "=SUM(C2:C30)"

# --- Match the saved selection state ---------------------------------------
$ws.Range("I2").Select()
